$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Add the two new worksheets at the end of the workbook
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$wsAnalysis = $wb.Worksheets.Add($null, $lastSheet)
$wsAnalysis.Name = "Analysis (Nothing)"

$wsWeather = $wb.Worksheets.Add($null, $wsAnalysis)
$wsWeather.Name = "Weather Data"

# ------------------------------------------------------------------
# 2. "Analysis (Nothing)" sheet contents
# ------------------------------------------------------------------
$wsAnalysis.Range("B2").Value = "This is the analysis for doing nothing"

$wsAnalysis.Range("B3").Value = "Analysis period"
$wsAnalysis.Range("C3").Value = 5
$wsAnalysis.Range("D3").Value = "years"

$wsAnalysis.Range("C4").Value = 20
$wsAnalysis.Range("D4").Value = "years"

$wsAnalysis.Range("C5").Value = 40
$wsAnalysis.Range("D5").Value = "years"

$wsAnalysis.Columns.Item(2).ColumnWidth = 33.86

# ------------------------------------------------------------------
# 3. "Weather Data" sheet contents
# ------------------------------------------------------------------
$wsWeather.Range("B2").Value = "Edmonton Weather Data"

$wsWeather.Range("B5").Value = "Longest Day"
$wsWeather.Range("C5").Value = 17
$wsWeather.Range("D5").Value = "hours"
$wsWeather.Range("F5").Value = "Sunrise"
$wsWeather.Range("G5").Value = 0.21111111111111111
$wsWeather.Range("G5").NumberFormat = "h:mm AM/PM"
$wsWeather.Range("I5").Value = "Sun highest"
$wsWeather.Range("J5").Value = 60
$wsWeather.Range("K5").Value = "deg"

$wsWeather.Range("B6").Value = "Shortest Day"
$wsWeather.Range("C6").Value = 7.6
$wsWeather.Range("D6").Value = "hours"
$wsWeather.Range("F6").Value = "Sunrise"
$wsWeather.Range("G6").Value = 0.36805555555555558
$wsWeather.Range("G6").NumberFormat = "h:mm AM/PM"
$wsWeather.Range("I6").Value = "Sun lowest"
$wsWeather.Range("J6").Value = 13
$wsWeather.Range("K6").Value = "deg"

$wsWeather.Range("A8:B8").Font.Bold = $true

$wsWeather.Columns.Item(2).ColumnWidth = 9.43
$wsWeather.Columns.Item(9).ColumnWidth = 11.29

$wsWeather.PageSetup.Orientation = 1

# ------------------------------------------------------------------
# 4. Selections / active sheet bookkeeping
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("I19").Select()

$wsAnalysis.Activate()
$wsAnalysis.Range("D17").Select()

$wsWeather.Activate()
$wsWeather.Range("B8").Select()
